$wb = $excel.ActiveWorkbook

# Layer0 sheet updates
$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.4053896379782324
$ws1.Range("C2").Value = -0.07217897947900613
$ws1.Range("B3").Value = 0.4936523938066099
$ws1.Range("C3").Value = -1.06621943263362
$ws1.Range("B4").Value = 0.7027033615272439
$ws1.Range("C4").Value = -0.4626479010554579

# Layer1 sheet updates
$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.8420369112401874
$ws2.Range("C2").Value = -0.5190313505020749
$ws2.Range("B3").Value = 0.5716933026895223
$ws2.Range("C3").Value = 0.9325260376499798
$ws2.Range("B4").Value = -1.555400529821802
$ws2.Range("C4").Value = 0.009879523655046487
